# Delete the entire "Group" column (column B) from the "nc" worksheet.
# This shifts Experiment/Category/numeric columns one position to the left,
# matching the target diff (A:G -> A:F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nc")
$ws.Activate()

$col = $ws.Columns.Item(2)
$null = $col.Select()
$null = $col.Delete()
